$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User Registration")

# ---------------------------------------------------------------------
# 1. Clear out the old "FR-REG-0X" / description cells so the now-unused
#    shared strings get pruned from the workbook's shared string table
#    when we save (matches the big block of <si> removals in the diff).
# ---------------------------------------------------------------------
$ws.Range("B8:C14").ClearContents()

# ---------------------------------------------------------------------
# 2. Re-populate the "User Registration" module table (rows 8-17) with
#    the new Test Scenario / Requirement / Description data.
# ---------------------------------------------------------------------
$ws.Range("A15").Value = "TS_REG_08"
$ws.Range("A16").Value = "TS_REG_09"
$ws.Range("A17").Value = "TS_REG_10"

$ws.Range("B12").Value = "FR_REG_05"
$ws.Range("B13").Value = "FR_REG_06"
$ws.Range("B14").Value = "FR_REG_07"

$ws.Range("B15").Value = "FR_REG_08"
$ws.Range("B16").Value = "FR_REG_09"
$ws.Range("B17").Value = "FR_REG_10"

$ws.Range("B9").Value = "FR_REG_02"
$ws.Range("B8").Value = "FR_REG_01"

$ws.Range("C8").Value = "Verify user can successfully register using valid email and password."

$ws.Range("B10").Value = "FR_REG_03"
$ws.Range("B11").Value = "FR_REG_04"

$ws.Range("C9").Value = "Verify system displays validation messages when mandatory registration fields are left blank."
$ws.Range("C10").Value = "Verify system displays appropriate error message for invalid email format during registration."
$ws.Range("C11").Value = "Verify system prevents registration using an already registered email address."
$ws.Range("C12").Value = "Verify system validates phone number format and length during registration."
$ws.Range("C13").Value = "Verify OTP is successfully sent to a valid phone number during registration."

$ws.Range("C17").Value = "Verify system allows user to resend OTP and generates a new OTP."
$ws.Range("C16").Value = "Verify system displays appropriate message when an expired OTP is entered."
$ws.Range("C15").Value = "Verify system displays error message when an invalid OTP is entered."

$ws.Range("C14").Value = "Verify user is able to complete registration by entering a valid OTP."

# ---------------------------------------------------------------------
# 3. Apply the new left/center/indent alignment style to the (now wider)
#    Description column for the data rows.
# ---------------------------------------------------------------------
$descRange = $ws.Range("C8:C17")
$descRange.HorizontalAlignment = -4131
$descRange.VerticalAlignment = -4108
$descRange.IndentLevel = 1

# ---------------------------------------------------------------------
# 4. Widen the Description column to fit the longer text.
# ---------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 88.59244791666667

# ---------------------------------------------------------------------
# 5. Update the view: select the cell below the new last row, and make
#    "User Registration" the active/selected sheet/tab.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("C19").Select()
